{"js": "// \"menambah latar belakang laporan\" \u2014 insert the background paragraph\n// right after the \"Latar belakang permasalahan\" heading (and before the\n// existing blank paragraph that precedes the \"Solusi\" heading).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst headingText = \"Latar belakang permasalahan\";\nlet heading = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === headingText) {\n    heading = paragraphs.items[i];\n    break;\n  }\n}\nif (!heading) {\n  throw new Error('Could not find heading paragraph \"' + headingText + '\"');\n}\n\n// The paragraph right after the heading is the existing empty spacer\n// paragraph before \"Solusi\". Insert the new body paragraph right before it\n// so it picks up the surrounding Normal formatting (no bold / no heading\n// style) instead of inheriting the heading's run/paragraph formatting.\nconst spacer = heading.getNextOrNullObject();\nspacer.load(\"text\");\nawait context.sync();\n\nconst newParagraphText =\n  \"Msalah kesehatan mental semakin menjadi perhatian utama di era modern ini. \" +\n  \"Tingkat stres, kecemasan, dan depresi semakin meningkat di kalangan masyarakat. \" +\n  \"Beban hidup sehari-hari, tuntutan pekerjaan, masalah pribadi, dan tekanan sosial \" +\n  \"dapat berdampak negatif pada kesejahteraan mental seseorang. Untuk mengatasi \" +\n  \"masalah ini, penting bagi individu untuk mendapatkan dukungan mental yang tepat \" +\n  \"dan aksesibilitas terhadap layanan kesehatan mental yang berkualitas.\";\n\nlet newPara;\nif (!spacer.isNullObject) {\n  newPara = spacer.insertParagraph(newParagraphText, \"Before\");\n} else {\n  // Fallback: heading was the last paragraph in the body.\n  newPara = heading.insertParagraph(newParagraphText, \"After\");\n  newPara.styleBuiltIn = Word.BuiltInStyleName.normal;\n}\n\n// 720 twips == 36 points -- Office.js paragraph indents are in points.\nnewPara.firstLineIndent = 36;\n\nawait context.sync();\n", "ps1": "# \"menambah latar belakang laporan\" \u2014 insert the background paragraph\n# right after the \"Latar belakang permasalahan\" heading (and before the\n# existing blank paragraph that precedes the \"Solusi\" heading).\n\n$d = $word.ActiveDocument\n\n$headingText = \"Latar belakang permasalahan\"\n$heading = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $headingText) {\n        $heading = $p\n        break\n    }\n}\nif ($null -eq $heading) {\n    throw \"Could not find heading paragraph '$headingText'\"\n}\n\n$newParagraphText = \"Msalah kesehatan mental semakin menjadi perhatian utama di era modern ini. \" +\n    \"Tingkat stres, kecemasan, dan depresi semakin meningkat di kalangan masyarakat. \" +\n    \"Beban hidup sehari-hari, tuntutan pekerjaan, masalah pribadi, dan tekanan sosial \" +\n    \"dapat berdampak negatif pada kesejahteraan mental seseorang. Untuk mengatasi \" +\n    \"masalah ini, penting bagi individu untuk mendapatkan dukungan mental yang tepat \" +\n    \"dan aksesibilitas terhadap layanan kesehatan mental yang berkualitas.\"\n\n# The paragraph right after the heading is the existing empty spacer\n# paragraph before \"Solusi\". Insert the new body paragraph right before it\n# so it picks up the surrounding Normal formatting (no bold / no heading\n# style) instead of inheriting the heading's run/paragraph formatting.\n$spacer = $heading.Next()\nif ($null -ne $spacer) {\n    $r = $spacer.Range\n    $r.Collapse(1)  # wdCollapseStart\n    $r.InsertParagraphBefore()\n    $newPara = $heading.Next()\n} else {\n    # Fallback: heading was the last paragraph in the body.\n    $r = $heading.Range\n    $r.Collapse(0)  # wdCollapseEnd\n    $r.InsertParagraphAfter()\n    $newPara = $heading.Next()\n    $newPara.Range.Style = $d.Styles.Item(\"Normal\")\n}\n\n$newPara.Range.Text = $newParagraphText\n$newPara.Format.FirstLineIndent = 36  # points (720 twips = 36pt = 0.5in)\n\nWrite-Output \"done\"\n"}
